$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("OrgData").Range("A2").Value = "AUTO_ORG_EMWRC"
$wb.Worksheets.Item("SubOrgData").Range("A2").Value = "AUTO_ORG_RBLQH"
$wb.Worksheets.Item("MemberData").Range("A2").Value = "First6563"
$wb.Worksheets.Item("StaffData").Range("A2").Value = "SFirstWSOPJ"
